$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "306.18"
Set-TextValue $ws.Range("E2") "0.92%"
Set-TextValue $ws.Range("D3") "36.18"
Set-TextValue $ws.Range("E3") "-1.47%"
Set-TextValue $ws.Range("D4") "5.059"
Set-TextValue $ws.Range("E4") "1.28%"
Set-TextValue $ws.Range("D5") "0.07937"
Set-TextValue $ws.Range("E5") "3.27%"
Set-TextValue $ws.Range("D6") "2.226"
Set-TextValue $ws.Range("E6") "7.85%"
Set-TextValue $ws.Range("E7") "0.69%"
Set-TextValue $ws.Range("D8") "4.140"
Set-TextValue $ws.Range("E8") "2.81%"
Set-TextValue $ws.Range("D9") "0.9277"
Set-TextValue $ws.Range("E9") "1.31%"
Set-TextValue $ws.Range("D10") "0.09811"
Set-TextValue $ws.Range("E10") "3.34%"
Set-TextValue $ws.Range("D11") "0.1872"
Set-TextValue $ws.Range("E11") "1.13%"
Set-TextValue $ws.Range("D12") "0.09159"
Set-TextValue $ws.Range("E12") "6.59%"
Set-TextValue $ws.Range("D13") "0.03707"
Set-TextValue $ws.Range("E13") "3.72%"
Set-TextValue $ws.Range("D14") "0.09917"
Set-TextValue $ws.Range("E14") "-0.66%"
Set-TextValue $ws.Range("D15") "0.001430"
Set-TextValue $ws.Range("E15") "-2.30%"
Set-TextValue $ws.Range("D16") "0.005637"
Set-TextValue $ws.Range("E16") "-1.38%"
Set-TextValue $ws.Range("D17") "3.461"
Set-TextValue $ws.Range("E17") "-0.01%"
Set-TextValue $ws.Range("E18") "18.67%"
Set-TextValue $ws.Range("E19") "-0.06%"
Set-TextValue $ws.Range("E20") "-1.11%"
Set-TextValue $ws.Range("D21") "5.090"
Set-TextValue $ws.Range("E21") "3.48%"
Set-TextValue $ws.Range("E22") "2.01%"
Set-TextValue $ws.Range("D23") "0.04557"
Set-TextValue $ws.Range("E23") "-0.82%"
Set-TextValue $ws.Range("D24") "0.001237"
Set-TextValue $ws.Range("E24") "0.18%"
Set-TextValue $ws.Range("D25") "0.004779"
Set-TextValue $ws.Range("E25") "-6.45%"
Set-TextValue $ws.Range("E26") "-7.29%"
Set-TextValue $ws.Range("E27") "74.20%"
Set-TextValue $ws.Range("D39") "0.01919"
Set-TextValue $ws.Range("E39") "10.15%"
Set-TextValue $ws.Range("D40") "0.04922"
Set-TextValue $ws.Range("E40") "6.38%"
Set-TextValue $ws.Range("D41") "0.007768"
Set-TextValue $ws.Range("E41") "1.43%"
Set-TextValue $ws.Range("D42") "0.1396"
Set-TextValue $ws.Range("E42") "0.37%"
Set-TextValue $ws.Range("D43") "0.007798"
Set-TextValue $ws.Range("E43") "0.68%"
Set-TextValue $ws.Range("D44") "0.002218"
Set-TextValue $ws.Range("E44") "2.59%"
Set-TextValue $ws.Range("D45") "0.01143"
Set-TextValue $ws.Range("E45") "10.00%"
Set-TextValue $ws.Range("D46") "0.00006282"
Set-TextValue $ws.Range("E46") "0.48%"
Set-TextValue $ws.Range("D47") "0.00000000749"
Set-TextValue $ws.Range("E47") "-0.18%"
Set-TextValue $ws.Range("D48") "52.39"
Set-TextValue $ws.Range("E48") "51.70%"
Set-TextValue $ws.Range("E49") "-10.14%"
Set-TextValue $ws.Range("D50") "0.00002098"
Set-TextValue $ws.Range("E50") "-0.18%"
Set-TextValue $ws.Range("E51") "-0.18%"
